$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Reference table gets a new "Precio" column (C) and one row's code changes ---
$ws.Range("C17").Value = 150
$ws.Range("C18").Value = 200
$ws.Range("A19").Value = 201
$ws.Range("C19").Value = 1200

# --- Main table: add lookups (B, C), quantities (D) and totals (E) for rows 12-14 ---
$ws.Range("B12").Formula = '=_xlfn.XLOOKUP(A12,A17:A19,B17:B19,"ERROR")'
$ws.Range("C12").Formula = '=_xlfn.XLOOKUP(A12,A17:A19,C17:C19,"ERROR")'
$ws.Range("D12").Value = 3
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Formula = "=D12*C12"
$ws.Range("E12").NumberFormat = "General"

$ws.Range("B13").Formula = '=_xlfn.XLOOKUP(A13,A18:A20,B18:B20,"ERROR")'
$ws.Range("C13").Formula = '=_xlfn.XLOOKUP(A13,A18:A20,C18:C20,"ERROR")'
$ws.Range("D13").Value = 4
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Formula = "=D13*C13"
$ws.Range("E13").NumberFormat = "General"

$ws.Range("B14").Formula = '=_xlfn.XLOOKUP(A14,A19:A21,B19:B21,"ERROR")'
$ws.Range("C14").Formula = '=_xlfn.XLOOKUP(A14,A19:A21,C19:C21,"ERROR")'
$ws.Range("D14").Value = 12
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Formula = "=D14*C14"
$ws.Range("E14").NumberFormat = "General"

# --- New row 15: a lookup code that isn't present, to show the ERROR branch ---
$ws.Range("A15").Value = 203
$ws.Range("B15").Formula = '=_xlfn.XLOOKUP(A15,A20:A22,B20:B22,"ERROR")'
$ws.Range("C15").Formula = '=_xlfn.XLOOKUP(A15,A20:A22,C20:C22,"ERROR")'

# Move the active selection
$ws.Range("G14").Select()

$wb.Save()
